# Auto-generated: updates Price (D) and Volume(1h) (E) columns
# for the crypto-symbol refresh performed by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. Values are written with a leading
# apostrophe so Excel keeps them as literal text (matching the
# source data, which stores prices/percentages as text, not numbers)
# instead of re-interpreting them as numbers/percentages.
$updates = @{
    "D2" = "309.61"
    "E2" = "-3.73%"
    "D3" = "48.33"
    "E3" = "-1.09%"
    "D4" = "5.145"
    "E4" = "-3.55%"
    "D5" = "0.07762"
    "E5" = "-4.19%"
    "D6" = "4.471"
    "E6" = "-2.30%"
    "D7" = "1.312"
    "E7" = "19.97%"
    "D8" = "1.564"
    "E8" = "-6.49%"
    "D9" = "0.1226"
    "E9" = "-7.32%"
    "D10" = "0.1929"
    "E10" = "-1.39%"
    "D11" = "0.04681"
    "E11" = "3.28%"
    "D12" = "0.09282"
    "E12" = "-2.59%"
    "D13" = "0.1048"
    "E13" = "0.11%"
    "D14" = "0.001261"
    "E14" = "-5.40%"
    "D15" = "0.04169"
    "E15" = "-3.04%"
    "D16" = "0.005836"
    "E16" = "0.42%"
    "D17" = "3.331"
    "E17" = "-1.48%"
    "D18" = "2.236"
    "E18" = "-8.00%"
    "E19" = "2.88%"
    "D20" = "8.260"
    "E20" = "1.02%"
    "D21" = "0.1358"
    "E21" = "-2.95%"
    "D22" = "0.3031"
    "E22" = "3.59%"
    "D23" = "0.001272"
    "E23" = "-2.62%"
    "D24" = "0.004132"
    "E24" = "-3.03%"
    "D25" = "0.0001350"
    "E25" = "0.19%"
    "E26" = "-3.96%"
    "D38" = "0.02586"
    "E38" = "-6.91%"
    "D39" = "0.05864"
    "E39" = "5.84%"
    "D40" = "0.01076"
    "E40" = "70.99%"
    "D41" = "0.007895"
    "E41" = "1.89%"
    "D42" = "0.1422"
    "E42" = "-1.63%"
    "D43" = "0.008422"
    "E43" = "9.67%"
    "D44" = "0.007655"
    "E44" = "-13.21%"
    "D45" = "0.3117"
    "E45" = "-11.54%"
    "D46" = "0.00006929"
    "E46" = "0.98%"
    "D47" = "0.00000000750"
    "E47" = "0.04%"
    "D48" = "0.05669"
    "E48" = "6.35%"
    "E49" = "0.24%"
    "D50" = "0.00002099"
    "E50" = "0.04%"
    "D51" = "0.0001999"
    "E51" = "0.04%"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
